# Update gh-pages output values on "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6652
$ws1.Range("F16").Value = 3317
$ws1.Range("F19").Value = 1960
$ws1.Range("F20").Value = 74

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6652
$ws4.Range("F17").Value = 3317
$ws4.Range("F20").Value = 1960
$ws4.Range("F21").Value = 74
